# "alteracoes na pagina usuario" - mark several functional requirements as
# Done on the "Funcionais" sheet, and fix up the header cell's leftover
# (no-op) style so it matches the style used for the same header on the
# other requirement sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Funcionais")

# The "Descrição" header cell (C2) carried a stray/duplicate style (an
# applyAlignment flag with no actual alignment). Re-apply bold (its font
# already is bold) so it re-resolves onto the canonical bold+border style.
$ws.Range("C2").Font.Bold = $true

# Flip the "Situação" column to "Done" for the requirements that have been
# completed.
$ws.Range("D3").Value = "Done"
$ws.Range("D8").Value = "Done"
$ws.Range("D9").Value = "Done"
$ws.Range("D10").Value = "Done"
$ws.Range("D11").Value = "Done"
$ws.Range("D12").Value = "Done"

$ws.Range("D3").Select()
